$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.691.74"
$ws.Range("E2").Value = "  +1.51%  "
$ws.Range("D3").Value = "1.877.69"
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.79%  "
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4755"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2823"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.64%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06490"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.63"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.18%  "
$ws.Range("D11").Value = "1.876.53"
$ws.Range("E11").Value = "  +2.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07579"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.15"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +14.06%  "
$ws.Range("E14").Value = "  +2.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6467"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.39%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +31.77%  "
$ws.Range("D17").Value = "30.686.35"
$ws.Range("E17").Value = "  +1.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.11%  "
$ws.Range("E19").Value = "  +5.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007496"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.95%  "
$ws.Range("D21").Value = "2.123.07"
$ws.Range("E21").Value = "  +2.50%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.128"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.93%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.110"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.42%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "169.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.213"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.941"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.14%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1055"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.36%  "
$ws.Range("E30").Value = "  -1.63%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.148"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.938"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.65%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04978"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.18%  "
$ws.Range("E34").Value = "  +2.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7163"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.66%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.711"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.83%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01902"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.700"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.040"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8941"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.74"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("E42").Value = "  -0.11%  "
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.560"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.29"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +8.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.302"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.71%  "
$ws.Range("E47").Value = "  +1.39%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.852"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.49%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.43"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.52%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05605"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.381"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.21%  "
